$wb = $excel.ActiveWorkbook

# --- Sheet "318.15": selection moved to the header row (A1:C1) ---
$ws3 = $wb.Worksheets.Item("318.15")
$ws3.Activate()
$ws3.Range("A1:C1").Select()

# --- Sheet "333.15": a new header row is inserted at the top, pushing the
#     existing 66 data rows down to rows 2-67, and the header cells get the
#     shared-string labels "x_1" / "x_2" / "x_3" (same labels already used
#     as headers on the other sheets). ---
$ws4 = $wb.Worksheets.Item("333.15")
$ws4.Activate()
$ws4.Rows("1:1").Insert()
$ws4.Range("A1").Value = "x_1"
$ws4.Range("B1").Value = "x_2"
$ws4.Range("C1").Value = "x_3"

# The sheet also carries a left-over (empty) sort-range marker (F1:I33)
# from an earlier "Data > Sort" operation. Inserting the row above it
# shifts that marker down by one row too (F2:I34), so reproduce that via
# a no-op sort over the same (empty) range.
$ws4.Sort.SortFields.Clear()
$ws4.Sort.SortFields.Add($ws4.Range("F2:F34"), 0, 2, 0, 0)
$ws4.Sort.SetRange($ws4.Range("F2:I34"))
$ws4.Sort.Header = 0
$ws4.Sort.Apply()

# Final active selection on the front sheet
$ws4.Range("E8").Select()
